# Daily attendance processing - 2025-10-08 13:50:17
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Row 3 (Year3/C1 ANATOMY session 2) - reorder "Recorded By" list
$ws.Range("G3").Value = "Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"

# Row 4 (Year3/C1 ANATOMY session 3) - reorder "Recorded By" list and update attendance count
$ws.Range("G4").Value = "Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("H4").Value = "105/221"

# Row 10 - Average Attendance % for ANATOMY/C1 group, recalculated after H4 update
# (leading ' keeps this a literal text value instead of Excel auto-converting it to a percentage number)
$ws.Range("L10").Value = "'36.5%"

# Row 12 (Year3/C1 HISTOLOGY session 1) - reorder "Recorded By" list
$ws.Range("G12").Value = "mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"

# Row 15 - Avg Attendance % for HISTOLOGY/C1 group
$ws.Range("S15").Value = "'41.8%"

# Row 16 - Avg Attendance % for PHARMACOLOGY/C1 group
$ws.Range("S16").Value = "'32.6%"

# Row 25 (Year3/C2 ANATOMY session 2) - reorder "Recorded By" list
$ws.Range("G25").Value = "Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"

# Row 26 (Year3/C2 ANATOMY session 3) - reorder "Recorded By" list and update attendance count
$ws.Range("G26").Value = "Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("H26").Value = "87/246"

# Row 34 (Year3/C2 HISTOLOGY session 1) - reorder "Recorded By" list
$ws.Range("G34").Value = "mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"

# Row 41 (Year3/C2 PHYSIOLOGY session 1) - reorder "Recorded By" list
$ws.Range("G41").Value = "Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
